$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 4303
$ws.Range("I64").Value = 3942.7144
$ws.Range("J64").Value = 4513.1665
$ws.Range("K64").Value = 3942.7144
$ws.Range("L64").Value = 4513.1665
$ws.Range("M64").Value = -3694.7144
$ws.Range("N64").Value = -5009.1665
$ws.Range("H67").Value = 4303
$ws.Range("I67").Value = 3942.7144
$ws.Range("J67").Value = 4513.1665
$ws.Range("K67").Value = 3942.7144
$ws.Range("L67").Value = 4513.1665
$ws.Range("M67").Value = -3084.7144
$ws.Range("N67").Value = -6229.1665
$ws.Range("H76").Value = 4677.722
$ws.Range("I76").Value = 3075.125
$ws.Range("J76").Value = 5959.8
$ws.Range("K76").Value = 3075.125
$ws.Range("L76").Value = 5959.8
$ws.Range("M76").Value = -2760.125
$ws.Range("N76").Value = -6589.8
$ws.Range("H79").Value = 4677.722
$ws.Range("I79").Value = 3075.125
$ws.Range("J79").Value = 5959.8
$ws.Range("K79").Value = 3075.125
$ws.Range("L79").Value = 5959.8
$ws.Range("M79").Value = -1983.125
$ws.Range("N79").Value = -8143.8
$ws.Range("H98").Value = 17781792
$ws.Range("I98").Value = 21167198
$ws.Range("J98").Value = 8411.25
$ws.Range("K98").Value = 21167198
$ws.Range("L98").Value = 8411.25
$ws.Range("M98").Value = -21165700
$ws.Range("N98").Value = -11407.25
$ws.Range("H122").Value = 17781792
$ws.Range("I122").Value = 21167198
$ws.Range("J122").Value = 8411.25
$ws.Range("K122").Value = 63501594
$ws.Range("L122").Value = 25233.75
$ws.Range("M122").Value = -63499144
$ws.Range("N122").Value = -30133.75
$ws.Range("H138").Value = 24192.305
$ws.Range("I138").Value = 11751.294
$ws.Range("J138").Value = 31485.31
$ws.Range("K138").Value = 35253.882
$ws.Range("L138").Value = 94455.93000000001
$ws.Range("M138").Value = -30113.882
$ws.Range("N138").Value = -104735.93

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 24169.252
$ws.Range("I32").Value = 21287.5
$ws.Range("J32").Value = 37391.41
$ws.Range("K32").Value = 21287.5
$ws.Range("L32").Value = 37391.41
$ws.Range("M32").Value = -21000.5
$ws.Range("N32").Value = -37965.41
$ws.Range("H61").Value = 4389185.5
$ws.Range("I61").Value = 5558766
$ws.Range("J61").Value = 3257.1875
$ws.Range("K61").Value = 5558766
$ws.Range("L61").Value = 3257.1875
$ws.Range("M61").Value = -5558554
$ws.Range("N61").Value = -3681.1875
$ws.Range("H74").Value = 1192.4333
$ws.Range("I74").Value = 782.5
$ws.Range("K74").Value = 782.5
$ws.Range("M74").Value = 91.5
$ws.Range("H77").Value = 1192.4333
$ws.Range("I77").Value = 782.5
$ws.Range("K77").Value = 3912.5
$ws.Range("M77").Value = 455.5
$ws.Range("H88").Value = 29383
$ws.Range("I88").Value = 1216.6666
$ws.Range("J88").Value = 46282.8
$ws.Range("K88").Value = 1216.6666
$ws.Range("L88").Value = 46282.8
$ws.Range("M88").Value = -810.6666
$ws.Range("N88").Value = -47094.8
$ws.Range("H91").Value = 29383
$ws.Range("I91").Value = 1216.6666
$ws.Range("J91").Value = 46282.8
$ws.Range("K91").Value = 1216.6666
$ws.Range("L91").Value = 46282.8
$ws.Range("M91").Value = 187.3334
$ws.Range("N91").Value = -49090.8
$ws.Range("H136").Value = 4389185.5
$ws.Range("I136").Value = 5558766
$ws.Range("J136").Value = 3257.1875
$ws.Range("K136").Value = 16676298
$ws.Range("L136").Value = 9771.5625
$ws.Range("M136").Value = -16673748
$ws.Range("N136").Value = -14871.5625
$ws.Range("H138").Value = 45109.832
$ws.Range("J138").Value = 45109.832
$ws.Range("L138").Value = 45109.832
$ws.Range("N138").Value = -55389.832

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H60").Value = 30000
$ws.Range("J60").Value = 30000
$ws.Range("L60").Value = 30000
$ws.Range("N60").Value = -31198
$ws.Range("H86").Value = 2215.7334
$ws.Range("I86").Value = 1580.5834
$ws.Range("K86").Value = 1580.5834
$ws.Range("M86").Value = -457.5834
$ws.Range("H89").Value = 2215.7334
$ws.Range("I89").Value = 1580.5834
$ws.Range("K89").Value = 7902.916999999999
$ws.Range("M89").Value = -2286.916999999999
$ws.Range("H105").Value = 3416.6667
$ws.Range("I105").Value = 4733.3335
$ws.Range("J105").Value = 2100
$ws.Range("K105").Value = 4733.3335
$ws.Range("L105").Value = 2100
$ws.Range("M105").Value = -2986.3335
$ws.Range("N105").Value = -5594
$ws.Range("H134").Value = 10041924
$ws.Range("I134").Value = 11914661
$ws.Range("J134").Value = 1974748
$ws.Range("K134").Value = 35743983
$ws.Range("L134").Value = 5924244
$ws.Range("M134").Value = -35741448
$ws.Range("N134").Value = -5929314

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 13347464
$ws.Range("I31").Value = 23817574
$ws.Range("J31").Value = 21868.908
$ws.Range("K31").Value = 23817574
$ws.Range("L31").Value = 21868.908
$ws.Range("M31").Value = -23817279
$ws.Range("N31").Value = -22458.908
$ws.Range("H34").Value = 13347464
$ws.Range("I34").Value = 23817574
$ws.Range("J34").Value = 21868.908
$ws.Range("K34").Value = 23817574
$ws.Range("L34").Value = 21868.908
$ws.Range("M34").Value = -23817372
$ws.Range("N34").Value = -22272.908
$ws.Range("H62").Value = 17005
$ws.Range("J62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("N62").ClearContents()
$ws.Range("H65").Value = 17005
$ws.Range("J65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("N65").ClearContents()

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H26").Value = 4630840
$ws.Range("I26").Value = 766.6667
$ws.Range("J26").Value = 6945876.5
$ws.Range("K26").Value = 2300.0001
$ws.Range("L26").Value = 20837629.5
$ws.Range("M26").Value = -2012.0001
$ws.Range("N26").Value = -20838205.5
$ws.Range("H29").Value = 97.72727
$ws.Range("J29").Value = 185.8
$ws.Range("L29").Value = 557.4000000000001
$ws.Range("N29").Value = -1111.4
$ws.Range("H39").Value = 2769.3333
$ws.Range("J39").Value = 2987.2
$ws.Range("L39").Value = 8961.599999999999
$ws.Range("N39").Value = -9549.599999999999

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 32086.834
$ws.Range("I70").Value = 59417.11
$ws.Range("J70").Value = 4756.5557
$ws.Range("K70").Value = 59417.11
$ws.Range("L70").Value = 4756.5557
$ws.Range("M70").Value = -59147.11
$ws.Range("N70").Value = -5296.5557
$ws.Range("H73").Value = 32086.834
$ws.Range("I73").Value = 59417.11
$ws.Range("J73").Value = 4756.5557
$ws.Range("K73").Value = 59417.11
$ws.Range("L73").Value = 4756.5557
$ws.Range("M73").Value = -58481.11
$ws.Range("N73").Value = -6628.5557
$ws.Range("H80").Value = 2933.3333
$ws.Range("J80").Value = 3020
$ws.Range("L80").Value = 3020
$ws.Range("N80").Value = -5016
$ws.Range("H83").Value = 2933.3333
$ws.Range("J83").Value = 3020
$ws.Range("L83").Value = 15100
$ws.Range("N83").Value = -25084

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 2156.9697
$ws.Range("I61").Value = 1573.6364
$ws.Range("J61").Value = 3323.6365
$ws.Range("K61").Value = 1573.6364
$ws.Range("L61").Value = 3323.6365
$ws.Range("M61").Value = -1371.6364
$ws.Range("N61").Value = -3727.6365
$ws.Range("H113").Value = 2156.9697
$ws.Range("I113").Value = 1573.6364
$ws.Range("J113").Value = 3323.6365
$ws.Range("K113").Value = 1573.6364
$ws.Range("L113").Value = 3323.6365
$ws.Range("M113").Value = 596.3635999999999
$ws.Range("N113").Value = -7663.636500000001

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 18455530
$ws.Range("I132").Value = 31282996
$ws.Range("J132").Value = 608617.7
$ws.Range("K132").Value = 93848988
$ws.Range("L132").Value = 1825853.1
$ws.Range("M132").Value = -93846458
$ws.Range("N132").Value = -1830913.1
